# Data Modeling.xlsx — build out migrations + seeding data and create initial DB
#
# The diff renames a couple of stale placeholder values in the ER-diagram /
# seed-data mock-up on Sheet1 and updates the saved window/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Charities table (rows 2-36): fix misspelled charity name ------------
$ws.Range("C18").Value = "Lutheran Services in America"

# --- Donations table (rows 2-36): FK column header rename ---------------
# "person_id (FK)" -> "user_id (FK)" now that the referenced table is Users
$ws.Range("H4").Value = "user_id (FK)"

# --- Users table (rows 43-49): seed data cleanup -------------------------
# Row 47 (was "Nathan Thomas" / "nate" / "nate"): split into real first/last
# name columns and give the account its actual username
$ws.Range("D47").Value = "Thomas"
$ws.Range("C47").Value = "Nathan"
$ws.Range("E47").Value = "nwthomas"

# Row 46 ("Admin" account): last_name column should read "Admin", not "admin"
$ws.Range("D46").Value = "Admin"

# --- Saved view state ------------------------------------------------------
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$win.Zoom = 140
$ws.Range("D46").Select() | Out-Null
